$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 4 and row 5 for columns D, J, K, L, M, N, O, P, Q
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

foreach ($col in $cols) {
    $addr4 = $col + "4"
    $addr5 = $col + "5"
    $val4 = $ws.Range($addr4).Value()
    $val5 = $ws.Range($addr5).Value()
    $ws.Range($addr4).Value = $val5
    $ws.Range($addr5).Value = $val4
}
